# Extend the year/value table from column K into a new column L,
# mirroring the existing 2020 figures (row 3: year label, row 4: value),
# then leave the selection where the user's cursor ended up (L10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column L, row 3 (year header "2020") -----------------------------
# Copy K3's formatting into L3, then set the value explicitly (avoids
# any floating point drift PasteSpecial's value-copy could introduce).
[void]$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L3").Value = 2020

# --- Column L, row 4 (figure "6.18") -----------------------------------
[void]$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L4").Value = 6.18

# Clear the marching-ants clipboard marquee left behind by Copy().
$excel.CutCopyMode = 0

# --- Restore the final selection recorded in the saved file -----------
[void]$ws.Range("L10").Select()
